$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.966.87"
$ws.Range("E2").Value = "  -0.23%  "

$ws.Range("D3").Value = "1.872.49"
$ws.Range("E3").Value = "  +0.43%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.39"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.20%  "

$ws.Range("E6").Value = "  +0.16%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5060"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.47%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3660"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.28%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07196"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.67%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8937"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.10%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.71"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.30%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.879.98"
$ws.Range("E12").Value = "  +0.90%  "

$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07529"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.66%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "94.68"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +5.94%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.239"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.34%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.19%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008533"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.37%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.23"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.17%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9998"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.10%  "

$ws.Range("D20").Value = "27.016.92"
$ws.Range("E20").Value = "  -0.19%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.025"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.12%  "

$ws.Range("D22").Value = "2.105.49"
$ws.Range("E22").Value = "  +0.64%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.39"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.10%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.419"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.71%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.22"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.04%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.791"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.73%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.90"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.40%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.076"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.30%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "113.37"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.50%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.702"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.58%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.688"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.58%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09147"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.07%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05126"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.23%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7511"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +2.97%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.981"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.95%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.158"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.26%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.228"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +6.26%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.563"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +3.55%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5654"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +6.66%  "

$ws.Range("E40").Value = "  -1.92%  "

$ws.Range("E41").Value = "  -0.17%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.606"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.91%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "115.54"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.44%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.515"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +2.88%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1474"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.36%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4732"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.31%  "

$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9998"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.18%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.09"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.13%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.565"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.11%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.87"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.91%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.17"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.08%  "
